# Update data for 11-02-2020
# - Corrects the "Total Cases" value for 2020-10-29 (row 96)
# - Appends four new days of data (rows 97-100): 10-30, 10-31, 11-01, 11-02

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction on existing last row (row 96, 2020-10-29) ---
$ws.Range("B96").Value = 7121

# --- Extend formatting from the last existing row down to the four new rows ---
$ws.Range("A96:D96").Copy() | Out-Null
$ws.Range("A97:D100").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- New rows: Date / Total Cases ---
$ws.Range("A97").Value = 44134   # 2020-10-30
$ws.Range("B97").Value = 7170

$ws.Range("A98").Value = 44135   # 2020-10-31
$ws.Range("B98").Value = 7205

$ws.Range("A99").Value = 44136   # 2020-11-01
$ws.Range("B99").Value = 7241

$ws.Range("A100").Value = 44137  # 2020-11-02
$ws.Range("B100").Value = 7288

# --- New Cases column (difference from prior day) ---
$ws.Range("C97").Formula = "=B97-B96"
$ws.Range("C98").Formula = "=B98-B97"
$ws.Range("C99").Formula = "=B99-B98"
$ws.Range("C100").Formula = "=B100-B99"

# --- 7 Day Average column ---
$ws.Range("D97").Formula = "=AVERAGE(C91:C97)"
$ws.Range("D98").Formula = "=AVERAGE(C92:C98)"
$ws.Range("D99").Formula = "=AVERAGE(C93:C99)"
$ws.Range("D100").Formula = "=AVERAGE(C94:C100)"

# --- Match the final on-screen selection/scroll position ---
$ws.Range("A99:A100").Select() | Out-Null
